# "Cambios tileset y arreglo carga de tiles peques: engarzado normal"
#
# The tile-character lookup sheet used several "placeholder" Unicode glyphs
# (ő, ȫ, δ, σ, ē, ľ, Ĭ, ų, ă) in column B that are no longer part of the
# active tileset. They are replaced with a fresh batch of glyphs
# (×, ©, ¬, °, ¿, ¼, ½, ¾, Á) that the new, smaller tile font actually uses.
#
# Simply overwriting the cell values is enough: when the workbook is saved
# the now-unreferenced shared strings are dropped from the table and the
# new glyphs are appended to it, which automatically renumbers every other
# shared-string reference in the sheet (columns A, B and E) exactly as
# happened in the authoritative edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order matters: this reproduces the exact order the new glyphs were
# appended to the shared string table in the target workbook.
$ws.Range("B2").Value2  = [char]0x00D7   # "ų" -> "×"
$ws.Range("B28").Value2 = [char]0x00A9   # "ő" -> "©"
$ws.Range("B30").Value2 = [char]0x00AC   # "ȫ" -> "¬"
$ws.Range("B38").Value2 = [char]0x00B0   # "δ" -> "°"
$ws.Range("B39").Value2 = [char]0x00BF   # "σ" -> "¿"
$ws.Range("B5").Value2  = [char]0x00BC   # "ē" -> "¼"
$ws.Range("B49").Value2 = [char]0x00BD   # "ľ" -> "½"
$ws.Range("B52").Value2 = [char]0x00BE   # "Ĭ" -> "¾"
$ws.Range("B61").Value2 = [char]0x00C1   # "ă" -> "Á"

# Scroll the view back to the top and select B4 instead of the previous
# B57 (which also clears the stale topLeftCell="A37" scroll position).
[void]$ws.Range("B4").Select()
